$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 6192
$ws.Range("D20").Value = 5581084
$ws.Range("E20").Value = 901.3378552971576
$ws.Range("F20").Value = 6.961478666436349
$ws.Range("H20").Value = 26.24964569874315
